$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.062.74'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.680.56'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Formula = '="215.89"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("E6").Value = '  -2.84%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("D9").Formula = '="21.41"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E9").Value = '  +5.86%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '1.917.65'
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").Value = '1.680.95'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Formula = '="0.534"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Formula = '="66.32"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '27.047.85'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("E18").Value = '  +2.06%  '
$ws.Range("D19").Formula = '="236.45"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("D23").Formula = '="9.27"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("D24").Formula = '="2.14"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("D25").Formula = '="147.17"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").Formula = '="7.28"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("E27").Value = '  +3.92%  '
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Formula = '="1.17"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("D33").Value = '1.544.88'
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +5.52%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").Formula = '="0.916"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("E40").Value = '  +6.92%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  +3.17%  '
$ws.Range("D43").Formula = '="5.52"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E43").Value = '  -3.28%  '
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").Value = '1.822.44'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Formula = '="90.53"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  +3.25%  '
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("D51").Formula = '="8.00"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E51").Value = '  +6.31%  '
